$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the final week's attendance (column O) for each attendee (rows 3-8)
$ws.Range("O3:O8").Value = 1

# Update the SUM formula in column B to include the new column O.
# B3 is an independent formula.
$ws.Range("B3").Formula = "=SUM(C3:O3)"
# B4:B8 form a shared-formula group (B4 is the master); set the whole
# range at once via FormulaR1C1 so the shared-formula grouping is preserved.
$ws.Range("B4:B8").FormulaR1C1 = "=SUM(RC[1]:RC[13])"

# Remove the now-unused column P (trailing formula cell + empty cells)
$ws.Range("P1:P11").Delete()

# Update the active cell selection to reflect where the editor left off
$ws.Range("C10").Select() | Out-Null
